$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add the new "metadatos" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "datos"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadatos"

# --- Fill in the metadata table ---
$ws2.Range("A1").Value = "Variables"
$ws2.Range("B1").Value = "Descripción"
$ws2.Range("C1").Value = "Fuente"
$ws2.Range("D1").Value = "Fecha_de_extracción"

$ws2.Range("A2").Value = "anno"
$ws2.Range("B2").Value = "Año"
$ws2.Range("C2").Value = "…"
$ws2.Range("D2").Value = 45722

$ws2.Range("A3").Value = "codmpio"
$ws2.Range("B3").Value = "Código del municipio"
$ws2.Range("C3").Value = "…"
$ws2.Range("D3").Value = 45722

$ws2.Range("A4").Value = "numerador"
$ws2.Range("B4").Value = "# de casos de violencia interpersonal contra niños, niñas y adolescentes"
$ws2.Range("C4").Value = "Instituto Nacional de Medicina Legal y Ciencias Forenses"
$ws2.Range("D4").Value = 45722

$ws2.Range("A5").Value = "denominador"
$ws2.Range("B5").Value = "Total niños/niñas/adolescentes"
$ws2.Range("C5").Value = "Departamento Administrativo Nacional de Estadística (DANE)"
$ws2.Range("D5").Value = 45722

$ws2.Range("A6").Value = "interpersonal"
$ws2.Range("B6").Value = "x 100,000"
$ws2.Range("C6").Value = "Elaboración Propia"
$ws2.Range("D6").Value = 45722

# --- Apply the (non-scheme) Calibri font to the whole used range, including the trailing blank row ---
$ws2.Range("A1:D7").Font.Name = "Calibri"

# --- Apply the date format to the date column ---
$ws2.Range("D2:D6").NumberFormat = "d-mmm-yy"

# --- Match the recorded selection/active-cell state ---
$ws2.Range("D2:D6").Select()
$ws1.Range("C4").Select()
$ws2.Select()
